$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort data range A1:B7 (with header) descending by column B (TotalCpmI)
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B2:B7"), 0, 2)
$sort.SetRange($ws.Range("A1:B7"))
$sort.Header = 1
$sort.Apply()

# Apply number format (2 decimal places) to sorted values
$ws.Range("B2:B7").NumberFormat = "0.00"

# View settings
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("D8").Select()
